$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing E-column values (payout/prize totals) per season row
$ws.Range("E3").Value = 769172
$ws.Range("E4").Value = 862283
$ws.Range("E5").Value = 1004603
$ws.Range("E7").Value = 793013
$ws.Range("E8").Value = 996516
$ws.Range("E9").Value = 837545
$ws.Range("E10").Value = 854593
$ws.Range("E11").Value = 928845

# Add new rows 12 and 13 for the two new seasons
# Copy the formatting of the existing "A" column season-index cells (bold,
# bordered, centered style) onto the new ones before writing their values.
$ws.Range("A2").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "M2_11 Mahakam 2020"
$ws.Range("C12").Value = 9706
$ws.Range("D12").Value = 10783
$ws.Range("E12").Value = 983150
$ws.Range("F12").Value = 10000
$ws.Range("G12").Value = 10090
$ws.Range("H12").Value = 10216

$ws.Range("A2").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "M2_12 Wild Hunt 2020"
$ws.Range("C13").Value = 9756
$ws.Range("D13").Value = 10724
$ws.Range("E13").Value = 1182983
$ws.Range("F13").Value = 10070
$ws.Range("G13").Value = 10172
$ws.Range("H13").Value = 10313
